$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ridge)
$ws.Range("B2").Value = 0.1080945746927301
$ws.Range("C2").Value = -0.09408604021379063
$ws.Range("D2").Value = 0.2021806149065207
$ws.Range("E2").Value = 30.9769341665602
$ws.Range("F2").Value = 34.30873996740541
$ws.Range("G2").Value = "{'solver': 'saga', 'alpha': 1000.0}"
$ws.Range("H2").Value = 2.28

# Row 3 (Lasso)
$ws.Range("C3").Value = -0.07421678239518825
$ws.Range("D3").Value = 0.07421678239518825
$ws.Range("F3").Value = 33.99577893322876

# Row 4 (ElasticNet)
$ws.Range("H4").Value = 2.97

# Row 5 (SVR)
$ws.Range("H5").Value = 2.15

# Row 6 (KNN Regressor)
$ws.Range("H6").Value = 1.52

# Row 7 (Decision Tree)
$ws.Range("H7").Value = 2.3

# Row 8 (PLSRegression)
$ws.Range("H8").Value = 0.53
